$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Lost participants due to FixAlign problems"
$ws.Range("E2").Value = 24
$ws.Range("E3").Value = 27
$ws.Range("E4").Value = 35
$ws.Range("E5").Value = 52
$ws.Range("E6").Value = 84

$ws.Columns.Item(5).ColumnWidth = 35.7

$ws.Range("E7").Select()
